# Update values in column E (on the active sheet) to match the new
# algorithm result values described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E8"   = 16.468
    "E10"  = 16.483
    "E12"  = 17.792
    "E18"  = 16.593
    "E25"  = 17.478
    "E37"  = 16.643
    "E55"  = 16.494
    "E68"  = 17.584
    "E77"  = 16.715
    "E78"  = 16.535
    "E79"  = 17.268
    "E80"  = 16.481
    "E81"  = 16.465
    "E82"  = 16.776
    "E84"  = 16.669
    "E101" = 16.723
    "E102" = 16.615
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
